# Scheduled-runner market price refresh for Anima_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) per leve row
# across the ALC, ARM, BSM, CRP, CUL, LTW and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1077.9286
$ws.Range("I17").Value = 603.9583
$ws.Range("J17").Value = 1433.4062
$ws.Range("K17").Value = 1811.8749
$ws.Range("L17").Value = 4300.2186
$ws.Range("M17").Value = -1643.8749
$ws.Range("N17").Value = -4636.2186

# Row 21
$ws.Range("H21").Value = 17928.143
$ws.Range("I21").Value = 12000
$ws.Range("J21").Value = 18916.166
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 18916.166
$ws.Range("M21").Value = -11532
$ws.Range("N21").Value = -19852.166

# Row 23
$ws.Range("H23").Value = 17928.143
$ws.Range("I23").Value = 12000
$ws.Range("J23").Value = 18916.166
$ws.Range("K23").Value = 12000
$ws.Range("L23").Value = 18916.166
$ws.Range("M23").Value = -11766
$ws.Range("N23").Value = -19384.166

# Row 40
$ws.Range("H40").Value = 1152.7368
$ws.Range("I40").Value = 1081.8182
$ws.Range("J40").Value = 1250.25
$ws.Range("K40").Value = 1081.8182
$ws.Range("L40").Value = 1250.25
$ws.Range("M40").Value = -906.8181999999999
$ws.Range("N40").Value = -1600.25

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 82
$ws.Range("H82").Value = 2042
$ws.Range("I82").Value = 2042
$ws.Range("K82").Value = 6126
$ws.Range("M82").Value = -5720

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 85
$ws.Range("H85").Value = 2042
$ws.Range("I85").Value = 2042
$ws.Range("K85").Value = 6126
$ws.Range("M85").Value = -4722

# Row 94
$ws.Range("H94").Value = 2096.6667
$ws.Range("I94").Value = 2096.6667
$ws.Range("K94").Value = 2096.6667
$ws.Range("M94").Value = -1645.6667

# Row 100
$ws.Range("H100").Value = 2459.182
$ws.Range("J100").Value = 2707.1428
$ws.Range("L100").Value = 2707.1428
$ws.Range("N100").Value = -3789.1428

# Row 116
$ws.Range("H116").Value = 8859.058999999999
$ws.Range("I116").Value = 11691.363
$ws.Range("J116").Value = 3666.5
$ws.Range("K116").Value = 11691.363
$ws.Range("L116").Value = 3666.5
$ws.Range("M116").Value = -8249.362999999999
$ws.Range("N116").Value = -10550.5

# Row 132
$ws.Range("H132").Value = 2536.3447
$ws.Range("I132").Value = 2197.875
$ws.Range("J132").Value = 4161
$ws.Range("K132").Value = 6593.625
$ws.Range("L132").Value = 12483
$ws.Range("M132").Value = -4063.625
$ws.Range("N132").Value = -17543

# Row 136
$ws.Range("H136").Value = 36206.668
$ws.Range("J136").Value = 36206.668
$ws.Range("L136").Value = 36206.668
$ws.Range("N136").Value = -46406.668

# Row 138
$ws.Range("H138").Value = 2191.6
$ws.Range("I138").Value = 843.4
$ws.Range("J138").Value = 2769.4
$ws.Range("K138").Value = 2530.2
$ws.Range("L138").Value = 8308.200000000001
$ws.Range("M138").Value = 2609.8
$ws.Range("N138").Value = -18588.2

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3162.5
$ws.Range("I45").Value = 3116.6667
$ws.Range("K45").Value = 3116.6667
$ws.Range("M45").Value = -2739.6667

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 332.9
$ws.Range("I22").Value = 285
$ws.Range("J22").Value = 524.5
$ws.Range("K22").Value = 285
$ws.Range("L22").Value = 524.5
$ws.Range("M22").Value = -112
$ws.Range("N22").Value = -870.5

# Row 107
$ws.Range("H107").Value = 49340.24
$ws.Range("I107").Value = 92463.17999999999
$ws.Range("J107").Value = 1905
$ws.Range("K107").Value = 92463.17999999999
$ws.Range("L107").Value = 1905
$ws.Range("M107").Value = -90543.17999999999
$ws.Range("N107").Value = -5745

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 7813800
$ws.Range("I107").Value = 20834200
$ws.Range("J107").Value = 1560
$ws.Range("K107").Value = 20834200
$ws.Range("L107").Value = 1560
$ws.Range("M107").Value = -20832280
$ws.Range("N107").Value = -5400

# Row 134
$ws.Range("H134").Value = 7580506
$ws.Range("I134").Value = 11369636
$ws.Range("J134").Value = 2245.6365
$ws.Range("K134").Value = 34108908
$ws.Range("L134").Value = 6736.9095
$ws.Range("M134").Value = -34106373
$ws.Range("N134").Value = -11806.9095

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1034.8889
$ws.Range("I5").Value = 402
$ws.Range("K5").Value = 1206
$ws.Range("M5").Value = -1094

# Row 68
$ws.Range("H68").Value = 835.8889
$ws.Range("I68").Value = 1251
$ws.Range("J68").Value = 717.2857
$ws.Range("K68").Value = 3753
$ws.Range("L68").Value = 2151.8571
$ws.Range("M68").Value = -2942
$ws.Range("N68").Value = -3773.8571

# Row 71
$ws.Range("H71").Value = 835.8889
$ws.Range("I71").Value = 1251
$ws.Range("J71").Value = 717.2857
$ws.Range("K71").Value = 11259
$ws.Range("L71").Value = 6455.571300000001
$ws.Range("M71").Value = -7203
$ws.Range("N71").Value = -14567.5713

# Row 76
$ws.Range("H76").Value = 3600

# Row 79
$ws.Range("H79").Value = 3600

# Row 113
$ws.Range("H113").Value = 531.4194
$ws.Range("J113").Value = 558.4
$ws.Range("L113").Value = 1675.2
$ws.Range("N113").Value = -6015.2

# Row 135
$ws.Range("H135").Value = 1034.8889
$ws.Range("I135").Value = 402
$ws.Range("K135").Value = 3618
$ws.Range("M135").Value = -1083

# Row 140
$ws.Range("H140").Value = 1629.3334
$ws.Range("I140").Value = 1388.5714
$ws.Range("K140").Value = 4165.7142
$ws.Range("M140").Value = 1014.2858

# Row 141
$ws.Range("H141").Value = 4939.0386
$ws.Range("I141").Value = 3385
$ws.Range("J141").Value = 6271.0713
$ws.Range("K141").Value = 10155
$ws.Range("L141").Value = 18813.2139
$ws.Range("M141").Value = -4975
$ws.Range("N141").Value = -29173.2139

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 4076.4443
$ws.Range("I122").Value = 3859.158
$ws.Range("J122").Value = 4592.5
$ws.Range("K122").Value = 11577.474
$ws.Range("L122").Value = 13777.5
$ws.Range("M122").Value = -9127.474
$ws.Range("N122").Value = -18677.5

$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 26498
$ws.Range("J82").Value = 26498
$ws.Range("L82").Value = 26498
$ws.Range("N82").Value = -27264

# Row 85
$ws.Range("H85").Value = 26498
$ws.Range("J85").Value = 26498
$ws.Range("L85").Value = 26498
$ws.Range("N85").Value = -29150

# Row 105
$ws.Range("H105").Value = 95000
$ws.Range("J105").Value = 95000
$ws.Range("L105").Value = 95000
$ws.Range("N105").Value = -101988

# Row 136
$ws.Range("H136").Value = 3804.125
$ws.Range("I136").Value = 3775.0908
$ws.Range("K136").Value = 11325.2724
$ws.Range("M136").Value = -8775.2724
